$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update League names (col A) and Multiplier values (col B) for rows 2-36 ---
$ws.Cells.Item(2, 1).Value = "Australia A-League Men"
$ws.Cells.Item(2, 2).Value = 0.805
$ws.Cells.Item(3, 1).Value = "Austria 2. Liga"
$ws.Cells.Item(3, 2).Value = 0.626
$ws.Cells.Item(4, 1).Value = "Belgium Challenger Pro League"
$ws.Cells.Item(4, 2).Value = 0.776
$ws.Cells.Item(5, 1).Value = "Bulgaria First League"
$ws.Cells.Item(5, 2).Value = 0.832
$ws.Cells.Item(6, 1).Value = "Croatia 1. HNL"
$ws.Cells.Item(6, 2).Value = 1.067
$ws.Cells.Item(7, 1).Value = "Czech First Tier"
$ws.Cells.Item(7, 2).Value = 1.0
$ws.Cells.Item(8, 1).Value = "Denmark 1st Division"
$ws.Cells.Item(8, 2).Value = 0.62
$ws.Cells.Item(9, 1).Value = "England League One"
$ws.Cells.Item(9, 2).Value = 0.973
$ws.Cells.Item(10, 1).Value = "England League Two"
$ws.Cells.Item(10, 2).Value = 0.773
$ws.Cells.Item(11, 1).Value = "England National League"
$ws.Cells.Item(11, 2).Value = 0.515
$ws.Cells.Item(12, 1).Value = "Finland Veikkausliiga"
$ws.Cells.Item(12, 2).Value = 0.755
$ws.Cells.Item(13, 1).Value = "France National 1"
$ws.Cells.Item(13, 2).Value = 0.793
$ws.Cells.Item(14, 1).Value = "Germany 3. Liga"
$ws.Cells.Item(14, 2).Value = 0.808
$ws.Cells.Item(15, 1).Value = "Hungary NB I"
$ws.Cells.Item(15, 2).Value = 0.975
$ws.Cells.Item(16, 1).Value = "Iceland Besta Deild"
$ws.Cells.Item(16, 2).Value = 0.725
$ws.Cells.Item(17, 1).Value = "Italy Serie C"
$ws.Cells.Item(17, 2).Value = 0.65
$ws.Cells.Item(18, 1).Value = "Netherlands Eerste Divisie"
$ws.Cells.Item(18, 2).Value = 0.692
$ws.Cells.Item(19, 1).Value = "Norway 1. Division"
$ws.Cells.Item(19, 2).Value = 0.717
$ws.Cells.Item(20, 1).Value = "Poland 1 Liga"
$ws.Cells.Item(20, 2).Value = 0.739
$ws.Cells.Item(21, 1).Value = "Poland Ekstraklasa"
$ws.Cells.Item(21, 2).Value = 1.089
$ws.Cells.Item(22, 1).Value = "Republic of Ireland Premier Division"
$ws.Cells.Item(22, 2).Value = 0.811
$ws.Cells.Item(23, 1).Value = "Romania Liga 1"
$ws.Cells.Item(23, 2).Value = 0.967
$ws.Cells.Item(24, 1).Value = "Scotland Championship"
$ws.Cells.Item(24, 2).Value = 0.586
$ws.Cells.Item(25, 1).Value = "Scotland Premiership"
$ws.Cells.Item(25, 2).Value = 1.0
$ws.Cells.Item(26, 1).Value = "Serbia Super Liga"
$ws.Cells.Item(26, 2).Value = 0.835
$ws.Cells.Item(27, 1).Value = "Slovakia 1. Liga"
$ws.Cells.Item(27, 2).Value = 0.888
$ws.Cells.Item(28, 1).Value = "Slovenia 1. Liga"
$ws.Cells.Item(28, 2).Value = 0.874
$ws.Cells.Item(29, 1).Value = "South Africa Premier Division"
$ws.Cells.Item(29, 2).Value = 0.861
$ws.Cells.Item(30, 1).Value = "Sweden Superettan"
$ws.Cells.Item(30, 2).Value = 0.717
$ws.Cells.Item(31, 1).Value = "Switzerland Challenge League"
$ws.Cells.Item(31, 2).Value = 0.706
$ws.Cells.Item(32, 1).Value = "Tunisia Ligue 1"
$ws.Cells.Item(32, 2).Value = 0.72
$ws.Cells.Item(33, 1).Value = "USA USL Championship"
$ws.Cells.Item(33, 2).Value = 0.612
$ws.Cells.Item(34, 1).Value = "Portugal Liga 3"
$ws.Cells.Item(34, 2).Value = 0.636
$ws.Cells.Item(35, 1).Value = "Portugal Liga Revelacao Sub 23"
$ws.Cells.Item(35, 2).Value = 0.636
$ws.Cells.Item(36, 1).Value = "Canada Premier League"
$ws.Cells.Item(36, 2).Value = 0.605

# --- Rows 37-54 no longer hold league data: drop col A entirely, blank col B ---
$ws.Range("A37:A54").Clear()
$ws.Range("B55").Copy()
$ws.Range("B37:B54").PasteSpecial(-4122)
$ws.Range("B37:B54").ClearContents()

# --- Shrink the trailing blank filler rows from 1000 down to 979 ---
$ws.Range("A980:B1000").Clear()
